$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 30271
$ws.Range("E2").Value = 588036487030
$ws.Range("F2").Value = 2830314328
$ws.Range("G2").Value = 0.25133

$ws.Range("D3").Value = 1869.2
$ws.Range("E3").Value = 224695169815
$ws.Range("F3").Value = 4691237029
$ws.Range("G3").Value = 0.55328

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 83382209358
$ws.Range("F4").Value = 11307100051
$ws.Range("G4").Value = -0.02426

$ws.Range("D5").Value = 234.76
$ws.Range("E5").Value = 36590426871
$ws.Range("F5").Value = 280690954
$ws.Range("G5").Value = -0.5875

$ws.Range("E6").Value = 27508838298
$ws.Range("F6").Value = 1791758569
$ws.Range("G6").Value = 0.00314

$ws.Range("D7").Value = 0.469274
$ws.Range("E7").Value = 24521869218
$ws.Range("F7").Value = 299284846
$ws.Range("G7").Value = 0.35015

$ws.Range("D8").Value = 1868.78
$ws.Range("E8").Value = 14229098867
$ws.Range("F8").Value = 12222143
$ws.Range("G8").Value = 0.57869

$ws.Range("D9").Value = 0.286559
$ws.Range("E9").Value = 10042580816
$ws.Range("F9").Value = 157542023
$ws.Range("G9").Value = 0.17597

$ws.Range("D10").Value = 0.065873
$ws.Range("E10").Value = 9228668745
$ws.Range("F10").Value = 156099951
$ws.Range("G10").Value = 0.78345

$ws.Range("B11").Value = "SOL"
$ws.Range("C11").Value = "Solana"
$ws.Range("D11").Value = 21.6
$ws.Range("E11").Value = 8689837815
$ws.Range("F11").Value = 436603738
$ws.Range("G11").Value = -0.94263

$ws.Range("B12").Value = "TRX"
$ws.Range("C12").Value = "TRON"
$ws.Range("D12").Value = 0.07945000000000001
$ws.Range("E12").Value = 7138337119
$ws.Range("F12").Value = 278732944
$ws.Range("G12").Value = 0.28291

$ws.Range("B13").Value = "LTC"
$ws.Range("C13").Value = "Litecoin"
$ws.Range("D13").Value = 96.31999999999999
$ws.Range("E13").Value = 7062696026
$ws.Range("F13").Value = 729982986
$ws.Range("G13").Value = -0.7025

$ws.Range("B14").Value = "MATIC"
$ws.Range("C14").Value = "Polygon"
$ws.Range("D14").Value = 0.695423
$ws.Range("E14").Value = 6480679156
$ws.Range("F14").Value = 157131027
$ws.Range("G14").Value = 2.11276

$ws.Range("B15").Value = "DOT"
$ws.Range("C15").Value = "Polkadot"
$ws.Range("D15").Value = 5.1
$ws.Range("E15").Value = 6386024860
$ws.Range("F15").Value = 68629444
$ws.Range("G15").Value = -1.17097

$ws.Range("D16").Value = 268.28
$ws.Range("E16").Value = 5216489459
$ws.Range("F16").Value = 370646872
$ws.Range("G16").Value = -0.46994

$ws.Range("B17").Value = "AVAX"
$ws.Range("C17").Value = "Avalanche"
$ws.Range("D17").Value = 14.08
$ws.Range("E17").Value = 4868564536
$ws.Range("F17").Value = 217089971
$ws.Range("G17").Value = 4.06899

$ws.Range("B18").Value = "WBTC"
$ws.Range("C18").Value = "Wrapped Bitcoin"
$ws.Range("D18").Value = 30266
$ws.Range("E18").Value = 4767058142
$ws.Range("F18").Value = 37987294
$ws.Range("G18").Value = 0.22669

$ws.Range("D19").Value = 0.00000775
$ws.Range("E19").Value = 4581348393
$ws.Range("F19").Value = 136193478
$ws.Range("G19").Value = 5.43671

$ws.Range("D20").Value = 0.999682
$ws.Range("E20").Value = 4295568989
$ws.Range("F20").Value = 56144781
$ws.Range("G20").Value = 0.01237

$ws.Range("B21").Value = "BUSD"
$ws.Range("C21").Value = "Binance USD"
$ws.Range("D21").Value = 0.999781
$ws.Range("E21").Value = 4045267938
$ws.Range("F21").Value = 2184349579
$ws.Range("G21").Value = 0.0764

$ws.Range("B22").Value = "UNI"
$ws.Range("C22").Value = "Uniswap"
$ws.Range("D22").Value = 5.26
$ws.Range("E22").Value = 3966710646
$ws.Range("F22").Value = 68892017
$ws.Range("G22").Value = -0.98258

$ws.Range("D23").Value = 3.83
$ws.Range("E23").Value = 3556172755
$ws.Range("F23").Value = 1548675
$ws.Range("G23").Value = -0.82997

$ws.Range("D24").Value = 6.2
$ws.Range("E24").Value = 3208289369
$ws.Range("F24").Value = 85341069
$ws.Range("G24").Value = 0.38062

$ws.Range("B25").Value = "XMR"
$ws.Range("C25").Value = "Monero"
$ws.Range("D25").Value = 167.43
$ws.Range("E25").Value = 3038666908
$ws.Range("F25").Value = 62787149
$ws.Range("G25").Value = 0.24568

$ws.Range("B26").Value = "TUSD"
$ws.Range("C26").Value = "TrueUSD"
$ws.Range("D26").Value = 0.9996969999999999
$ws.Range("E26").Value = 2944752923
$ws.Range("F26").Value = 939878320
$ws.Range("G26").Value = 0.00602

$ws.Range("B27").Value = "ATOM"
$ws.Range("C27").Value = "Cosmos Hub"
$ws.Range("D27").Value = 9.359999999999999
$ws.Range("E27").Value = 2739496377
$ws.Range("F27").Value = 64654227
$ws.Range("G27").Value = 1.49827

$ws.Range("B28").Value = "XLM"
$ws.Range("C28").Value = "Stellar"
$ws.Range("D28").Value = 0.098839
$ws.Range("E28").Value = 2673991879
$ws.Range("F28").Value = 42479530
$ws.Range("G28").Value = 0.39792

$ws.Range("B29").Value = "ETC"
$ws.Range("C29").Value = "Ethereum Classic"
$ws.Range("D29").Value = 18.84
$ws.Range("E29").Value = 2672936077
$ws.Range("F29").Value = 76444042
$ws.Range("G29").Value = -0.1113

$ws.Range("D30").Value = 42.88
$ws.Range("E30").Value = 2572674305
$ws.Range("F30").Value = 1209900
$ws.Range("G30").Value = 0.62358

$ws.Range("D31").Value = 1.36
$ws.Range("E31").Value = 2012621476
$ws.Range("F31").Value = 6785348
$ws.Range("G31").Value = -1.30378

$ws.Range("B32").Value = "FIL"
$ws.Range("C32").Value = "Filecoin"
$ws.Range("D32").Value = 4.33
$ws.Range("E32").Value = 1880279223
$ws.Range("F32").Value = 64236209
$ws.Range("G32").Value = -0.66376

$ws.Range("D33").Value = 4.06
$ws.Range("E33").Value = 1776295710
$ws.Range("F33").Value = 10474983
$ws.Range("G33").Value = 0.28856

$ws.Range("B34").Value = "LDO"
$ws.Range("C34").Value = "Lido DAO"
$ws.Range("D34").Value = 1.95
$ws.Range("E34").Value = 1716789854
$ws.Range("F34").Value = 14448684
$ws.Range("G34").Value = -0.61533

$ws.Range("D35").Value = 0.04710601
$ws.Range("E35").Value = 1520906159
$ws.Range("F35").Value = 11022242
$ws.Range("G35").Value = 0.08827

$ws.Range("D36").Value = 102.49
$ws.Range("E36").Value = 1491014062
$ws.Range("F36").Value = 10162331
$ws.Range("G36").Value = -0.81063

$ws.Range("D37").Value = 7.1
$ws.Range("E37").Value = 1489407950
$ws.Range("F37").Value = 53610404
$ws.Range("G37").Value = -0.92474

$ws.Range("D38").Value = 0.056804
$ws.Range("E38").Value = 1485486803
$ws.Range("F38").Value = 3786181
$ws.Range("G38").Value = 0.62022

$ws.Range("B39").Value = "ARB"
$ws.Range("C39").Value = "Arbitrum"
$ws.Range("D39").Value = 1.13
$ws.Range("E39").Value = 1445281418
$ws.Range("F39").Value = 103634211
$ws.Range("G39").Value = 0.3552

$ws.Range("B40").Value = "VET"
$ws.Range("C40").Value = "VeChain"
$ws.Range("D40").Value = 0.01868981
$ws.Range("E40").Value = 1359078423
$ws.Range("F40").Value = 31871732
$ws.Range("G40").Value = -0.34258

$ws.Range("D41").Value = 1.33
$ws.Range("E41").Value = 1241815875
$ws.Range("F41").Value = 44290185
$ws.Range("G41").Value = -3.09645

$ws.Range("B42").Value = "GRT"
$ws.Range("C42").Value = "The Graph"
$ws.Range("D42").Value = 0.114577
$ws.Range("E42").Value = 1043481205
$ws.Range("F42").Value = 18369847
$ws.Range("G42").Value = -2.47137

$ws.Range("B43").Value = "AAVE"
$ws.Range("C43").Value = "Aave"
$ws.Range("D43").Value = 71.7
$ws.Range("E43").Value = 1036823892
$ws.Range("F43").Value = 58130553
$ws.Range("G43").Value = -5.10367

$ws.Range("B44").Value = "FRAX"
$ws.Range("C44").Value = "Frax"
$ws.Range("D44").Value = 0.99892
$ws.Range("E44").Value = 1003190654
$ws.Range("F44").Value = 4275913
$ws.Range("G44").Value = 0.04946

$ws.Range("B45").Value = "RETH"
$ws.Range("C45").Value = "Rocket Pool ETH"
$ws.Range("D45").Value = 2012.76
$ws.Range("E45").Value = 940280631
$ws.Range("F45").Value = 698388
$ws.Range("G45").Value = 0.55771

$ws.Range("B46").Value = "STX"
$ws.Range("C46").Value = "Stacks"
$ws.Range("D46").Value = 0.64361
$ws.Range("E46").Value = 895000675
$ws.Range("F46").Value = 8849550
$ws.Range("G46").Value = -0.43981

$ws.Range("B47").Value = "EGLD"
$ws.Range("C47").Value = "MultiversX"
$ws.Range("D47").Value = 34.49
$ws.Range("E47").Value = 884951931
$ws.Range("F47").Value = 6261926
$ws.Range("G47").Value = 1.2078

$ws.Range("B48").Value = "ALGO"
$ws.Range("C48").Value = "Algorand"
$ws.Range("D48").Value = 0.111579
$ws.Range("E48").Value = 845229390
$ws.Range("F48").Value = 22908520
$ws.Range("G48").Value = -0.05091

$ws.Range("B49").Value = "MKR"
$ws.Range("C49").Value = "Maker"
$ws.Range("D49").Value = 915.73
$ws.Range("E49").Value = 825393175
$ws.Range("F49").Value = 35168653
$ws.Range("G49").Value = -4.52713

$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "EOS"
$ws.Range("D50").Value = 0.732386
$ws.Range("E50").Value = 814539108
$ws.Range("F50").Value = 93094383
$ws.Range("G50").Value = 1.83906

$ws.Range("B51").Value = "BSV"
$ws.Range("C51").Value = "Bitcoin SV"
$ws.Range("D51").Value = 40.74
$ws.Range("E51").Value = 784863474
$ws.Range("F51").Value = 17431619
$ws.Range("G51").Value = -1.87153

Write-Output "Updated cryptocurrency data for 2023-07-09"